$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update midRate (F) and effectiveDate (H) values for rows 2-20 ---
$ws.Range("F2").Value = 0.034387103185724
$ws.Range("H2").Value = "'2021-10-28"
$ws.Range("F3").Value = 0.034386712974107
$ws.Range("H3").Value = "'2021-10-28"
$ws.Range("F4").Value = 0.034386334595541
$ws.Range("H4").Value = "'2021-10-28"
$ws.Range("F5").Value = 0.034385956225302
$ws.Range("H5").Value = "'2021-10-28"
$ws.Range("F6").Value = 0.03438557786339
$ws.Range("H6").Value = "'2021-10-28"
$ws.Range("F7").Value = 2.822387175072676
$ws.Range("H7").Value = "'2021-10-28"
$ws.Range("F8").Value = 0.311526479750779
$ws.Range("H8").Value = "'2021-10-28"
$ws.Range("F9").Value = 0.034384064499004
$ws.Range("H9").Value = "'2021-10-28"
$ws.Range("F10").Value = 0.034384442827612
$ws.Range("H10").Value = "'2021-10-28"
$ws.Range("F11").Value = 0.034383307866763
$ws.Range("H11").Value = "'2021-10-28"
$ws.Range("F12").Value = 0.034382917741276
$ws.Range("H12").Value = "'2021-10-28"
$ws.Range("F13").Value = 0.034383686178721
$ws.Range("H13").Value = "'2021-10-28"
$ws.Range("F14").Value = 0.034382539446228
$ws.Range("H14").Value = "'2021-10-28"
$ws.Range("F15").Value = 0.345794618744143
$ws.Range("H15").Value = "'2021-10-28"
$ws.Range("F16").Value = 0.034381782881104
$ws.Range("H16").Value = "'2021-10-28"
$ws.Range("F17").Value = 0.034381404611028
$ws.Range("H17").Value = "'2021-10-28"
$ws.Range("F18").Value = 0.034381026349275
$ws.Range("H18").Value = "'2021-10-28"
$ws.Range("F19").Value = 0.034380648095845
$ws.Range("H19").Value = "'2021-10-28"
$ws.Range("F20").Value = 2.822785524755829
$ws.Range("H20").Value = "'2021-10-28"

# --- Update selection to reflect the latest active cell ---
$ws.Range("H30").Select()

# --- Re-touch columns L:M so their column metadata is recorded explicitly
#     (splits the L:M / N:XFD column-format run while keeping identical
#     width/style) ---
$ws.Range($ws.Cells.Item(1,12), $ws.Cells.Item(1,13)).EntireColumn.Hidden = $false
